$d = $word.ActiveDocument

# 1. ERP: teste -> ERP: SAP ECC/4HANA
$d.Content.Find.Execute("ERP: teste", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ERP: SAP ECC/4HANA", 2)

# 2. Processos desejados: Adiantamentos no Cartão -> Processos desejados: (trailing space kept)
$d.Content.Find.Execute("Processos desejados: Adiantamentos no Cartão", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Processos desejados: ", 2)

# 3. Informações necessárias pelo ERP: Empresa -> Informações necessárias pelo ERP: (trailing space kept)
$d.Content.Find.Execute("Informações necessárias pelo ERP: Empresa", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Informações necessárias pelo ERP: ", 2)

# 4. Replace the big analysis paragraph (the one right after the "Análise Funcional Recomendada:" heading)
$br = [char]11
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt.StartsWith("Para realizar a integra")) {
        $r = $p.Range
        $r.End = $r.End - 1
        $newText = "A melhor análise funcional para realizar a integração do ERP SAP ECC/4HANA com o SaaS Paytrack seria a seguinte:" + $br + $br + `
            "1. Levantamento dos processos desejados de integração, como adiantamento, prestação de contas, entre outros." + $br + `
            "2. Identificação das informações necessárias pelo ERP repassadas pelo cliente, como bukrs (empresa), entre outros campos específicos do ERP." + $br + `
            "3. Mapeamento dos campos necessários para a integração, apresentando em formato de tabela com as nomenclaturas do ERP, como bukrs, e gerando um JSON de exemplo com os dados formatados de acordo com o ERP." + $br + `
            "4. Definição da comunicação Síncrona com os Webservices do cliente para a integração." + $br + `
            "5. Estabelecimento da ativação da Paytrack nas integrações, com o cliente disponibilizando um Webservice para consumir." + $br + `
            "6. Separação da análise funcional por cenário selecionado, com uma análise para adiantamento, uma para prestação de contas, entre outros cenários específicos." + $br + $br + `
            "Com essa abordagem detalhada e organizada, será possível realizar uma integração eficiente e eficaz entre o ERP SAP ECC/4HANA e o SaaS Paytrack, garantindo a troca de informações de forma precisa e adequada para atender às necessidades do cliente."
        $r.Text = $newText
        break
    }
}

# 5. Remove the "Empresa" / "Descrição para Empresa" row from the field-mapping table
$t = $d.Tables(1)
for ($i = $t.Rows.Count; $i -ge 1; $i--) {
    $row = $t.Rows($i)
    if ($row.Cells(1).Range.Text.StartsWith("Empresa")) {
        $row.Delete()
    }
}
